$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 16

$ws.Range("C3").Select()
